$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header cells for row 1
$ws.Range("E1").Value = "SEVERE"
$ws.Range("G1").Value = "Description"
$ws.Range("H1").Value = "Hint"

# Row 4 gets a new Severity entry
$ws.Range("E4").Value = "CRITICAL"

# New row 5 with a new testcase
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "f(x,y)=sin(x)+cos(y)"
$ws.Range("D5").Value = "Linux"
$ws.Range("E5").Value = "CRITICAL"
$ws.Range("F5").Value = "FAILED"
$ws.Range("G5").Value = "f(x,y)=sin((x+cos(y)))"
$ws.Range("H5").Value = "analyse the output tree"

# Column width adjustments (values chosen so the engine's internal
# pixel-rounding reproduces the target stored widths as closely as possible)
$ws.Columns.Item(2).ColumnWidth = 17.451822916666668
$ws.Columns.Item(7).ColumnWidth = 18.877604166666668
$ws.Columns.Item(8).ColumnWidth = 15.307291666666666

# Update selection to match final state
$ws.Range("H6").Select()
